# Update "ficha tecnica" worksheet: refresh names of current office holders
# and tidy up the sheet selection/formatting, mirroring a manual edit made
# in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Secretária de Estado de Planejamento e Gestão
$ws.Range("A4").Value = "Silvia Caroline Listgarten Dias"

# Secretário-Adjunto de Estado de Planejamento e Gestão
$ws.Range("A5").Value = "Rodrigo Ferreira Matias"

# Subsecretário de Planejamento e Orçamento (Assessoria de Inteligência de Dados)
$ws.Range("A7").Value = "Gabriel Braico Dornas"
$ws.Range("A7").WrapText = $true

# Reset selection/view back to a normal cell instead of the stray
# full-column selection left over from the previous edit.
$ws.Range("B15").Select() | Out-Null
